# Scheduled market-data refresh: update computed price/profit columns
# (currentAveragePrice[NQ/HQ], LevePrice[NQ/HQ], LeveProfit[NQ/HQ]) for the
# leves whose Universalis price snapshot changed, across all eight job
# sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 646.7442
$ws.Range("J17").Value = 646.7442
$ws.Range("L17").Value = 1940.2326
$ws.Range("N17").Value = -2276.2326
$ws.Range("H112").Value = 1999.75
$ws.Range("I112").Value = 1999
$ws.Range("J112").Value = 2000
$ws.Range("K112").Value = 5997
$ws.Range("L112").Value = 6000
$ws.Range("M112").Value = -4889
$ws.Range("N112").Value = -8216
$ws.Range("H134").Value = 98971.42999999999
$ws.Range("J134").Value = 98971.42999999999
$ws.Range("L134").Value = 98971.42999999999
$ws.Range("N134").Value = -109111.43
$ws.Range("H136").Value = 77977.14
$ws.Range("J136").Value = 77977.14
$ws.Range("L136").Value = 77977.14
$ws.Range("N136").Value = -88177.14
$ws.Range("H138").Value = 2183.8
$ws.Range("J138").Value = 2499
$ws.Range("L138").Value = 7497
$ws.Range("N138").Value = -17777
$ws.Range("H139").Value = 74130.89999999999
$ws.Range("J139").Value = 74130.89999999999
$ws.Range("L139").Value = 74130.89999999999
$ws.Range("N139").Value = -84410.89999999999
$ws.Range("H140").Value = 66637.875
$ws.Range("J140").Value = 71982.336
$ws.Range("L140").Value = 71982.336
$ws.Range("N140").Value = -82342.336

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H7").Value = 83468.336
$ws.Range("J7").Value = 83468.336
$ws.Range("L7").Value = 83468.336
$ws.Range("N7").Value = -83696.336
$ws.Range("H52").Value = 50947.223
$ws.Range("J52").Value = 50947.223
$ws.Range("L52").Value = 50947.223
$ws.Range("N52").Value = -51583.223
$ws.Range("H61").Value = 114974.11
$ws.Range("I61").Value = 4431.143
$ws.Range("K61").Value = 4431.143
$ws.Range("M61").Value = -4219.143
$ws.Range("H63").Value = 1994
$ws.Range("I63").Value = 1994
$ws.Range("K63").Value = 1994
$ws.Range("M63").Value = -1308
$ws.Range("H66").Value = 1994
$ws.Range("I66").Value = 1994
$ws.Range("K66").Value = 9970
$ws.Range("M66").Value = -6538
$ws.Range("H74").Value = 81953.38
$ws.Range("I74").Value = 169416
$ws.Range("K74").Value = 169416
$ws.Range("M74").Value = -168542
$ws.Range("H77").Value = 81953.38
$ws.Range("I77").Value = 169416
$ws.Range("K77").Value = 847080
$ws.Range("M77").Value = -842712
$ws.Range("H118").Value = 59966.855
$ws.Range("J118").Value = 59966.855
$ws.Range("L118").Value = 59966.855
$ws.Range("N118").Value = -63280.855
$ws.Range("H127").Value = 91397.14
$ws.Range("J127").Value = 91397.14
$ws.Range("L127").Value = 91397.14
$ws.Range("N127").Value = -101317.14
$ws.Range("H136").Value = 114974.11
$ws.Range("I136").Value = 4431.143
$ws.Range("K136").Value = 13293.429
$ws.Range("M136").Value = -10743.429

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H50").Value = 70505.5
$ws.Range("J50").Value = 74464.8
$ws.Range("L50").Value = 74464.8
$ws.Range("N50").Value = -75612.8
$ws.Range("H51").Value = 56747.2
$ws.Range("J51").Value = 56747.2
$ws.Range("L51").Value = 56747.2
$ws.Range("N51").Value = -57729.2
$ws.Range("H53").Value = 32731.2
$ws.Range("J53").Value = 32731.2
$ws.Range("L53").Value = 32731.2
$ws.Range("N53").Value = -33879.2
$ws.Range("H55").Value = 30414.5
$ws.Range("J55").Value = 30414.5
$ws.Range("L55").Value = 30414.5
$ws.Range("N55").Value = -30960.5
$ws.Range("H119").Value = 97673
$ws.Range("J119").Value = 97673
$ws.Range("L119").Value = 97673
$ws.Range("N119").Value = -107349
$ws.Range("H138").Value = 76664.44500000001
$ws.Range("J138").Value = 76664.44500000001
$ws.Range("L138").Value = 76664.44500000001
$ws.Range("N138").Value = -86944.44500000001
$ws.Range("H140").Value = 43498.934
$ws.Range("J140").Value = 43498.934
$ws.Range("L140").Value = 43498.934
$ws.Range("N140").Value = -53858.934

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H25").Value = 8000
$ws.Range("J25").Value = 8000
$ws.Range("L25").Value = 8000
$ws.Range("N25").Value = -8348
$ws.Range("H118").Value = 64685
$ws.Range("J118").Value = 64685
$ws.Range("L118").Value = 64685
$ws.Range("N118").Value = -67999
$ws.Range("H132").Value = 3078622
$ws.Range("I132").Value = 2527433.5
$ws.Range("K132").Value = 7582300.5
$ws.Range("M132").Value = -7579770.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 113197.555
$ws.Range("J68").Value = 127034.75
$ws.Range("L68").Value = 381104.25
$ws.Range("N68").Value = -382726.25
$ws.Range("H71").Value = 113197.555
$ws.Range("J71").Value = 127034.75
$ws.Range("L71").Value = 1143312.75
$ws.Range("N71").Value = -1151424.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H13").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("L13").Value = 0
# LeveProfitHQ no longer computable for this leve (no HQ price data) - clear it
$ws.Range("N13").ClearContents()
$ws.Range("H109").Value = 25340.666
$ws.Range("J109").Value = 27570.75
$ws.Range("L109").Value = 27570.75
$ws.Range("N109").Value = -29650.75
$ws.Range("H116").Value = 56343.4
$ws.Range("J116").Value = 58890.5
$ws.Range("L116").Value = 58890.5
$ws.Range("N116").Value = -68068.5
$ws.Range("H119").Value = 80790.71000000001
$ws.Range("J119").Value = 80790.71000000001
$ws.Range("L119").Value = 80790.71000000001
$ws.Range("N119").Value = -90466.71000000001
$ws.Range("H126").Value = 2981.7222
$ws.Range("I126").Value = 2442
$ws.Range("J126").Value = 4385
$ws.Range("K126").Value = 7326
$ws.Range("L126").Value = 13155
$ws.Range("M126").Value = -4856
$ws.Range("N126").Value = -18095
$ws.Range("H140").Value = 94552.37
$ws.Range("J140").Value = 94957.60000000001
$ws.Range("L140").Value = 94957.60000000001
$ws.Range("N140").Value = -105317.6
$ws.Range("H141").Value = 72000
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H112").Value = 0
$ws.Range("J112").Value = 0
$ws.Range("L112").Value = 0
# LeveProfitHQ no longer computable for this leve (no HQ price data) - clear it
$ws.Range("N112").ClearContents()
$ws.Range("H117").Value = 84096.664
$ws.Range("J117").Value = 84096.664
$ws.Range("L117").Value = 84096.664
$ws.Range("N117").Value = -93274.664
$ws.Range("H118").Value = 85401.11
$ws.Range("J118").Value = 85401.11
$ws.Range("L118").Value = 85401.11
$ws.Range("N118").Value = -88715.11
$ws.Range("H132").Value = 4397.727
$ws.Range("I132").Value = 4437.5
$ws.Range("J132").Value = 4000
$ws.Range("K132").Value = 13312.5
$ws.Range("L132").Value = 12000
$ws.Range("M132").Value = -10782.5
$ws.Range("N132").Value = -17060

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H17").Value = 5000125
$ws.Range("I17").Value = 5000125
$ws.Range("K17").Value = 5000125
$ws.Range("M17").Value = -4999953
$ws.Range("H39").Value = 28044
$ws.Range("I39").Value = 28044
$ws.Range("K39").Value = 28044
$ws.Range("M39").Value = -27631
$ws.Range("H126").Value = 4682.8184
$ws.Range("I126").Value = 4017.375
$ws.Range("J126").Value = 6457.3335
$ws.Range("K126").Value = 12052.125
$ws.Range("L126").Value = 19372.0005
$ws.Range("M126").Value = -9582.125
$ws.Range("N126").Value = -24312.0005
$ws.Range("H127").Value = 60390
$ws.Range("I127").Value = 60390
$ws.Range("K127").Value = 60390
$ws.Range("M127").Value = -55430
$ws.Range("H132").Value = 2163.5417
$ws.Range("I132").Value = 1762.5555
$ws.Range("J132").Value = 3366.5
$ws.Range("K132").Value = 5287.666499999999
$ws.Range("L132").Value = 10099.5
$ws.Range("M132").Value = -2757.666499999999
$ws.Range("N132").Value = -15159.5
